$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-6)
$ws.Range("A2").Value = 4
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 1

# Update column B values (rows 2-3); B4:B6 unchanged
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4

# Remove rows 7 and 8 (delete entire rows so dimension shrinks to A1:B6)
$ws.Rows("7:8").Delete()
